# Update Write Latency min (O) and max (P) columns per updated test-5 generator output
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = "'1301"
$ws.Range("P3").Value = "1637.5k"

$ws.Range("O4").Value = "'1362"
$ws.Range("P4").Value = "1705.3k"

$ws.Range("O5").Value = "'1261"
$ws.Range("P5").Value = "6556.3k"

$ws.Range("O6").Value = "'1383"
$ws.Range("P6").Value = "1848.1k"

$ws.Range("O7").Value = "'1386"
$ws.Range("P7").Value = "2496.2k"

$ws.Range("O8").Value = "'1765"
$ws.Range("P8").Value = "12929k"

$ws.Range("O9").Value = "'1392"
$ws.Range("P9").Value = "8060.7k"

$ws.Range("O10").Value = "'1418"
$ws.Range("P10").Value = "2415.0k"

$ws.Range("O11").Value = "'1393"
$ws.Range("P11").Value = "'620273"

$ws.Range("O12").Value = "'2"
$ws.Range("P12").Value = "'514"

$ws.Range("O13").Value = "'1372"
$ws.Range("P13").Value = "'688717"

$ws.Range("P14").Value = "'1487"

$ws.Range("O15").Value = "'1567"
$ws.Range("P15").Value = "8542.8k"

$ws.Range("O16").Value = "'1443"
$ws.Range("P16").Value = "5451.8k"

$ws.Range("O17").Value = "'1594"
$ws.Range("P17").Value = "1519.0k"

$ws.Range("O18").Value = "'1230"
$ws.Range("P18").Value = "6902.5k"

$ws.Range("O19").Value = "'1383"
$ws.Range("P19").Value = "5156.5k"

$ws.Range("O20").Value = "'2"
$ws.Range("P20").Value = "'1429"

$ws.Range("O21").Value = "'2"
$ws.Range("P21").Value = "'6784"

$ws.Range("O22").Value = "'1262"
$ws.Range("P22").Value = "2159.5k"

$ws.Range("O23").Value = "'1393"
$ws.Range("P23").Value = "1582.5k"
